# Apply updated cryptocurrency price/volume data to Sheet1
# (values sourced from coinranking.com snapshot refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.185.27'
$ws.Range('E2').Value = '  +0.57%  '
# Row 3
$ws.Range('D3').Value = '3.314.02'
$ws.Range('E3').Value = '  +0.80%  '
# Row 4
$s = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = $s
$ws.Range('E4').Value = '  -0.03%  '
# Row 5
$s = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.53'
$ws.Range('D5').Style = $s
$ws.Range('E5').Value = '  +1.23%  '
# Row 6
$s = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.85'
$ws.Range('D6').Style = $s
$ws.Range('E6').Value = '  +1.04%  '
# Row 7
$s = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = $s
$ws.Range('E7').Value = '  -0.09%  '
# Row 8
$ws.Range('D8').Value = '3.308.82'
$ws.Range('E8').Value = '  +0.85%  '
# Row 9
$s = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.574'
$ws.Range('D9').Style = $s
$ws.Range('E9').Value = '  -2.44%  '
# Row 10
$s = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.176'
$ws.Range('D10').Style = $s
$ws.Range('E10').Value = '  -4.94%  '
# Row 11
$s = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.573'
$ws.Range('D11').Style = $s
$ws.Range('E11').Value = '  -1.85%  '
# Row 12
$s = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.06'
$ws.Range('D12').Style = $s
$ws.Range('E12').Value = '  -2.73%  '
# Row 13
$ws.Range('E13').Value = '  -0.90%  '
# Row 14
$ws.Range('D14').Value = '3.845.68'
$ws.Range('E14').Value = '  +1.13%  '
# Row 15
$s = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.45'
$ws.Range('D15').Style = $s
$ws.Range('E15').Value = '  -2.29%  '
# Row 16
$s = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '592.83'
$ws.Range('D16').Style = $s
$ws.Range('E16').Value = '  -7.90%  '
# Row 17
$ws.Range('D17').Value = '66.038.54'
$ws.Range('E17').Value = '  +0.38%  '
# Row 19
$ws.Range('D19').Value = '3.318.43'
$ws.Range('E19').Value = '  +1.05%  '
# Row 20
$ws.Range('E20').Value = '  -2.25%  '
# Row 21
$s = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.92'
$ws.Range('D21').Style = $s
$ws.Range('E21').Value = '  -3.85%  '
# Row 22
$s = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.894'
$ws.Range('D22').Style = $s
$ws.Range('E22').Value = '  -1.07%  '
# Row 23
$s = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.06'
$ws.Range('D23').Style = $s
$ws.Range('E23').Value = '  -1.81%  '
# Row 24
$s = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.02'
$ws.Range('D24').Style = $s
$ws.Range('E24').Value = '  +2.07%  '
# Row 25
$s = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.16'
$ws.Range('D25').Style = $s
# Row 26
$s = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.98'
$ws.Range('D26').Style = $s
$ws.Range('E26').Value = '  +0.22%  '
# Row 27
$s = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.71'
$ws.Range('D27').Style = $s
$ws.Range('E27').Value = '  +0.70%  '
# Row 28
$s = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.39'
$ws.Range('D28').Style = $s
$ws.Range('E28').Value = '  -1.75%  '
# Row 29
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$s = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.44'
$ws.Range('D29').Style = $s
$ws.Range('E29').Value = '  -2.65%  '
# Row 30
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$s = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.59'
$ws.Range('D30').Style = $s
$ws.Range('E30').Value = '  +0.93%  '
# Row 31
$s = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.63'
$ws.Range('D31').Style = $s
$ws.Range('E31').Value = '  +5.41%  '
# Row 32
$s = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.69'
$ws.Range('D32').Style = $s
$ws.Range('E32').Value = '  -6.23%  '
# Row 33
$s = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '561.15'
$ws.Range('D33').Style = $s
$ws.Range('E33').Value = '  +8.16%  '
# Row 34
$s = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.82'
$ws.Range('D34').Style = $s
$ws.Range('E34').Value = '  -2.06%  '
# Row 35
$ws.Range('D35').Value = '3.784.56'
$ws.Range('E35').Value = '  -0.04%  '
# Row 36
$ws.Range('E36').Value = '  -1.37%  '
# Row 37
$s = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = $s
$ws.Range('E37').Value = '  -0.03%  '
# Row 38
$s = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.90'
$ws.Range('D38').Style = $s
$ws.Range('E38').Value = '  -2.87%  '
# Row 39
$s = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.27'
$ws.Range('D39').Style = $s
$ws.Range('E39').Value = '  +0.94%  '
# Row 40
$ws.Range('E40').Value = '  -2.37%  '
# Row 41
$ws.Range('D41').Value = '0.0₃0686'
$ws.Range('E41').Value = '  -6.58%  '
# Row 42
$s = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.14'
$ws.Range('D42').Style = $s
$ws.Range('E42').Value = '  -7.43%  '
# Row 43
$s = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.37'
$ws.Range('D43').Style = $s
$ws.Range('E43').Value = '  +4.53%  '
# Row 44
$s = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.58'
$ws.Range('D44').Style = $s
$ws.Range('E44').Value = '  -5.22%  '
# Row 45
$s = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.332'
$ws.Range('D45').Style = $s
$ws.Range('E45').Value = '  -1.36%  '
# Row 46
$s = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0410'
$ws.Range('D46').Style = $s
$ws.Range('E46').Value = '  -1.02%  '
# Row 47
$ws.Range('E47').Value = '  -8.97%  '
# Row 48
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$s = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.126'
$ws.Range('D48').Style = $s
$ws.Range('E48').Value = '  -2.43%  '
# Row 49
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$s = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').Style = $s
$ws.Range('E49').Value = '  +0.08%  '
# Row 50
$s = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.52'
$ws.Range('D50').Style = $s
$ws.Range('E50').Value = '  -3.07%  '
# Row 51
$s = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '128.41'
$ws.Range('D51').Style = $s
$ws.Range('E51').Value = '  +5.37%  '
